$wb = $excel.ActiveWorkbook

# --- About sheet ---
$wsAbout = $wb.Worksheets.Item("About")

# Unit labels: "trillion ..." -> "million ..." (note source typo "milliom" preserved)
$wsAbout.Range("A12").Value = "milliom passenger-miles"
$wsAbout.Range("A13").Value = "million freight ton-miles"

# Conversion note label: "miles to kilometers" -> "km to miles"
$wsAbout.Range("A15").Value = "km to miles"

# Fix the conversion factor: was a hardcoded miles-to-km constant (1.60934),
# now computed as the reciprocal (km-to-miles)
$wsAbout.Range("B15").Formula = "=1/1.60934"

# --- CDCF-PMpPDOU sheet ---
$wsPM = $wb.Worksheets.Item("CDCF-PMpPDOU")
$wsPM.Range("B2").Formula = "=10^6*About!`$B`$15"

# --- CDCF-FTMpFDOU sheet ---
$wsFTM = $wb.Worksheets.Item("CDCF-FTMpFDOU")
$wsFTM.Range("B2").Formula = "=10^6*About!`$B`$15"

# --- View / selection state ---
# Set selections on the non-active sheets first (selecting a range on a sheet
# implicitly activates it), then activate + select on "About" last so it
# ends up as the active/visible tab.
$wsPM.Range("B4").Select() | Out-Null

$wsFTM.Range("B3").Select() | Out-Null

$wsAbout.Activate() | Out-Null
$wsAbout.Range("B16").Select() | Out-Null
